# Apply updated crypto price/volume figures to Sheet1 (rows 2-51).
# For D-column values that look numeric (e.g. "1.027"), force the cell
# to remain a text value (matching the source data, which stores these
# as text) by briefly switching the cell to a text NumberFormat before
# assigning the value, then resetting the style back to Normal so no
# stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.388.65'
$ws.Range("E2").Value = '  +3.69%  '
$ws.Range("D3").Value = '1.837.51'
$ws.Range("E3").Value = '  +3.62%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '1.027'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  +2.24%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '317.79'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.61%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '1.024'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.99%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.4360'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.26%  '
$ws.Range("E8").Value = '  +2.77%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.07344'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +2.97%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.8706'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +3.84%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '21.31'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +4.57%  '
$ws.Range("D12").Value = '1.915.60'
$ws.Range("E12").Value = '  +6.25%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '5.464'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +4.12%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.683'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +3.71%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.07102'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +2.79%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '82.19'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +3.98%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '1.029'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +2.09%  '
$ws.Range("E18").Value = '  +3.25%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '1.023'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.96%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '15.37'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +2.91%  '
$ws.Range("D21").Value = '27.414.35'
$ws.Range("E21").Value = '  +3.69%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '5.245'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("E23").Value = '  +0.75%  '
$ws.Range("D24").Value = '2.136.93'
$ws.Range("E24").Value = '  +5.70%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '156.72'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +2.91%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '1.905'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +5.43%  '
$ws.Range("E27").Value = '  +2.94%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '5.237'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +3.42%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.919'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +8.02%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '115.55'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +1.42%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.09033'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +1.71%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.197'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +7.40%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.7584'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.55%  '
$ws.Range("E34").Value = '  +3.27%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '2.862'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +4.13%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '1.025'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +2.11%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.146'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +4.06%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.01954'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +3.84%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.05243'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.05%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.5155'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +5.09%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '2.789'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +7.57%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.1660'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.90%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '6.544'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +3.08%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '8.463'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +6.27%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '108.35'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +3.45%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '10.54'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +3.81%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.026'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +2.26%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.674'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +2.66%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.4617'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +4.22%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.06299'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +1.96%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '1.876'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +9.30%  '

Write-Host "Applied 95 cell updates across rows 2-51"
